# Generate Report for Handoff
# Adds a new localization-status row for file
# "f759cc0e-c025-4549-b079-85e853f13c84.md" to the Overview, zh-cn and
# de-de sheets (mirrors the existing "Ready for handoff" rows already
# present on row 8 of each sheet).

$wb = $excel.ActiveWorkbook

$guid        = "f759cc0e-c025-4549-b079-85e853f13c84"
$mdName      = "$guid.md"
$mdPath      = "e2e\$guid.md"
$zhXlf       = "$guid.0e76342da4f5c3b031d1dbd7c11f1b277553e53a.zh-cn.xlf"
$deXlf       = "$guid.0e76342da4f5c3b031d1dbd7c11f1b277553e53a.de-de.xlf"
$hoDate      = "2016-09-05 20:51:37"
$zhHoDate    = "2016-09-05 20:51:32"
$deHoDate    = "2016-09-05 20:51:37"

# ---------------------------------------------------------------------
# Sheet "Overview" -> new row 9
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$rowOv = 9
$wsOverview.Cells.Item($rowOv, 1).Value = $mdName
$wsOverview.Cells.Item($rowOv, 3).Value = ".md"
$wsOverview.Cells.Item($rowOv, 4).Value = ""
$wsOverview.Cells.Item($rowOv, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item($rowOv, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item($rowOv, 7).Value = $hoDate
$wsOverview.Cells.Item($rowOv, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Cells.Item($rowOv, 2),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c946a384fb83ad5e89ad10b5e4d9318e1a6d0c75/e2e/$guid.md",
    [Type]::Missing,
    [Type]::Missing,
    $mdPath
) | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G9"))

# ---------------------------------------------------------------------
# Sheet "zh-cn" -> new row 9
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$rowZh = 9
$wsZhCn.Cells.Item($rowZh, 2).Value  = ".md"
$wsZhCn.Cells.Item($rowZh, 3).Value  = "Ready for handoff"
$wsZhCn.Cells.Item($rowZh, 4).Value  = "e2e"
$wsZhCn.Cells.Item($rowZh, 5).Value  = "ht"
$wsZhCn.Cells.Item($rowZh, 6).Value  = "False"
$wsZhCn.Cells.Item($rowZh, 7).Value  = $zhXlf
$wsZhCn.Cells.Item($rowZh, 8).Value  = $zhHoDate
$wsZhCn.Cells.Item($rowZh, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item($rowZh, 9).Value  = ""
$wsZhCn.Cells.Item($rowZh, 10).Value = ""
$wsZhCn.Cells.Item($rowZh, 11).Value = "0001-01-01 00:00:00"
$wsZhCn.Cells.Item($rowZh, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item($rowZh, 12).Value = ""
$wsZhCn.Cells.Item($rowZh, 13).Value = "True"
$wsZhCn.Cells.Item($rowZh, 14).Value = ""
$wsZhCn.Cells.Item($rowZh, 15).Value = "False"
$wsZhCn.Cells.Item($rowZh, 16).Value = ""

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Cells.Item($rowZh, 1),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c946a384fb83ad5e89ad10b5e4d9318e1a6d0c75/e2e/$guid.md",
    [Type]::Missing,
    [Type]::Missing,
    $mdName
) | Out-Null

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P9"))

# ---------------------------------------------------------------------
# Sheet "de-de" -> new row 9
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$rowDe = 9
$wsDeDe.Cells.Item($rowDe, 2).Value  = ".md"
$wsDeDe.Cells.Item($rowDe, 3).Value  = "Ready for handoff"
$wsDeDe.Cells.Item($rowDe, 4).Value  = "e2e"
$wsDeDe.Cells.Item($rowDe, 5).Value  = "ht"
$wsDeDe.Cells.Item($rowDe, 6).Value  = "False"
$wsDeDe.Cells.Item($rowDe, 7).Value  = $deXlf
$wsDeDe.Cells.Item($rowDe, 8).Value  = $deHoDate
$wsDeDe.Cells.Item($rowDe, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item($rowDe, 9).Value  = ""
$wsDeDe.Cells.Item($rowDe, 10).Value = ""
$wsDeDe.Cells.Item($rowDe, 11).Value = "0001-01-01 00:00:00"
$wsDeDe.Cells.Item($rowDe, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item($rowDe, 12).Value = ""
$wsDeDe.Cells.Item($rowDe, 13).Value = "True"
$wsDeDe.Cells.Item($rowDe, 14).Value = ""
$wsDeDe.Cells.Item($rowDe, 15).Value = "False"
$wsDeDe.Cells.Item($rowDe, 16).Value = ""

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Cells.Item($rowDe, 1),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c946a384fb83ad5e89ad10b5e4d9318e1a6d0c75/e2e/$guid.md",
    [Type]::Missing,
    [Type]::Missing,
    $mdName
) | Out-Null

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P9"))
